# Sanchez_Yehoshua_ProblemSolving.docx -- "Identified the patterns and how
# they helped me come to the solutions."
#
# The old closing of the document was two paragraphs: one holding only the
# `_GoBack` bookmark, and a final empty centered paragraph. This edit
# replaces that pair with three floating tables (the "First/Ring" finger
# counting charts for 0-100, 100-200 and 0-1000), a run of blank spacer
# paragraphs underneath them, and a closing paragraph that carries the new
# write-up plus the original `_GoBack` bookmark.

$d = $word.ActiveDocument

# Anchor on the `_GoBack` bookmark rather than a hard-coded paragraph index
# so the edit still lands correctly even if earlier content shifts.
$bm = $d.Bookmarks.Item("_GoBack")
$anchorPara = $bm.Range.Paragraphs.Item(1)
$lastPara = $d.Paragraphs.Last

# Remove the bookmark paragraph through the trailing empty centered
# paragraph (this collapses them down to a single empty paragraph mark,
# since a document always needs a final paragraph mark).
$killRange = $d.Range($anchorPara.Range.Start, $lastPara.Range.End)
$killRange.Delete()

# Insert the three tables + spacer paragraphs + new closing paragraph
# (carrying the write-up text and the restored `_GoBack` bookmark) at that
# now-collapsed point, right before the section break.
$insertPoint = $d.Paragraphs.Last.Range
$insertPoint.Collapse(0)
$insertPoint.InsertXML('<w:tbl xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:tblPr><w:tblStyle w:val="TableGrid"/><w:tblpPr w:leftFromText="180" w:rightFromText="180" w:vertAnchor="page" w:horzAnchor="page" w:tblpX="1369" w:tblpY="1981"/><w:tblW w:w="0" w:type="auto"/><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="696"/><w:gridCol w:w="685"/></w:tblGrid><w:tr><w:trPr><w:trHeight w:val="673"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="696" w:type="dxa"/></w:tcPr><w:p/></w:tc><w:tc><w:tcPr><w:tcW w:w="685" w:type="dxa"/></w:tcPr><w:p/></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="625"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="696" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>First</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="685" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>Ring</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="673"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="696" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>0</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="685" w:type="dxa"/></w:tcPr><w:p/></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="673"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="696" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>10</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="685" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="625"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="696" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="685" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>20</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="673"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="696" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="685" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>30</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="673"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="696" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>40</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="685" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="673"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="696" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>50</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="685" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="673"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="696" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="685" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>60</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="673"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="696" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="685" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>70</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="673"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="696" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>80</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="685" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="673"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="696" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>90</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="685" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="489"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="696" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="685" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>100</w:t></w:r></w:p></w:tc></w:tr></w:tbl><w:tbl><w:tblPr><w:tblStyle w:val="TableGrid"/><w:tblpPr w:leftFromText="180" w:rightFromText="180" w:vertAnchor="text" w:horzAnchor="page" w:tblpX="4249" w:tblpY="-130"/><w:tblW w:w="0" w:type="auto"/><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="851"/><w:gridCol w:w="797"/></w:tblGrid><w:tr><w:trPr><w:trHeight w:val="659"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="851" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>First</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="797" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>Ring</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="659"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="851" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="797" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>100</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="613"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="851" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="797" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>110</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="659"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="851" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>120</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="797" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="659"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="851" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>130</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="797" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="659"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="851" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="797" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>140</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="659"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="851" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="797" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>150</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="659"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="851" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>160</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="797" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="659"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="851" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>170</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="797" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="659"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="851" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="797" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>180</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="659"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="851" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="797" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>190</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="659"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="851" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>200</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="797" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p></w:tc></w:tr></w:tbl><w:tbl><w:tblPr><w:tblStyle w:val="TableGrid"/><w:tblpPr w:leftFromText="180" w:rightFromText="180" w:vertAnchor="text" w:horzAnchor="page" w:tblpX="7669" w:tblpY="-145"/><w:tblW w:w="1449" w:type="dxa"/><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="748"/><w:gridCol w:w="701"/></w:tblGrid><w:tr><w:trPr><w:trHeight w:val="662"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="748" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>First</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="701" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>Ring</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="607"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="748" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>0</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="701" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="662"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="748" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="701" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>100</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="662"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="748" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>200</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="701" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="607"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="748" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="701" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>300</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="662"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="748" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>400</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="701" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="662"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="748" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="701" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>500</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="662"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="748" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>600</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="701" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="662"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="748" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="701" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>700</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="662"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="748" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>800</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="701" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="662"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="748" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="701" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>900</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:trHeight w:val="662"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="748" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>1000</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="701" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p></w:tc></w:tr></w:tbl><w:p/><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p><w:p><w:r><w:t>With the last three graphs we use the established patterns to figure out the solutions. In the first we use the pattern to find that when counting the girl will land on her ring finger when hitting 100. Then continuing on we can see that every hundred the pattern alternates between the two fingers. So we take that pattern into the third table and see that the girl would land on 1,000 when counting on her fingers.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>')
